# ---------------------------------------------------------------------------
# Scheduled-runner market data refresh for the Chocobo Profits workbook.
#
# This mirrors the periodic bot that re-pulls Universalis current-average-price
# market data (columns H:N -- currentAveragePrice / currentAveragePriceNQ /
# currentAveragePriceHQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ /
# LeveProfitHQ) for specific Leve rows across the job sheets (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR) and writes the refreshed numbers back as static
# values (the sheets store plain numbers, not formulas).
#
# Where a profit column has no meaningful value for a row (e.g. the HQ price
# data dried up) the cell is cleared entirely rather than zeroed, matching
# how the bot omits the cell from its feed.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 26: Everything Is Impossible
$ws.Range("H26").Value = 42000
$ws.Range("J26").Value = 42000
$ws.Range("L26").Value = 42000
$ws.Range("N26").Value = -42688

# Row 33: Glazed and Confused
$ws.Range("H33").Value = 220.73334
$ws.Range("I33").Value = 137.36363
$ws.Range("J33").Value = 450
$ws.Range("K33").Value = 137.36363
$ws.Range("L33").Value = 450
$ws.Range("M33").Value = 91.63637
$ws.Range("N33").Value = -908

# Row 44: Alive and Unwell
$ws.Range("H44").Value = 20272.727
$ws.Range("J44").Value = 20272.727
$ws.Range("L44").Value = 20272.727
$ws.Range("N44").Value = -21196.727

# Row 52: Your Courtesy Wake-up Call
$ws.Range("H52").Value = 4000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 4000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 12000
$ws.Range("N52").Value = -12320
$ws.Range("M52").ClearContents()

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 3180.1177
$ws.Range("I137").Value = 1927.8462
$ws.Range("K137").Value = 5783.5386
$ws.Range("M137").Value = -3233.5386

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 20351
$ws.Range("I141").Value = 37000.332
$ws.Range("J141").Value = 3701.6667
$ws.Range("K141").Value = 111000.996
$ws.Range("L141").Value = 11105.0001
$ws.Range("M141").Value = -105820.996
$ws.Range("N141").Value = -21465.0001

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 7123.58
$ws.Range("I32").Value = 4531.643
$ws.Range("J32").Value = 10422.409
$ws.Range("K32").Value = 4531.643
$ws.Range("L32").Value = 10422.409
$ws.Range("M32").Value = -4244.643
$ws.Range("N32").Value = -10996.409

# Row 88: The Mast Chance
$ws.Range("H88").Value = 5130425.5
$ws.Range("I88").Value = 8335343
$ws.Range("J88").Value = 2557.4
$ws.Range("K88").Value = 8335343
$ws.Range("L88").Value = 2557.4
$ws.Range("M88").Value = -8334937
$ws.Range("N88").Value = -3369.4

# Row 91: The Rose and the Riveter (L)
$ws.Range("H91").Value = 5130425.5
$ws.Range("I91").Value = 8335343
$ws.Range("J91").Value = 2557.4
$ws.Range("K91").Value = 8335343
$ws.Range("L91").Value = 2557.4
$ws.Range("M91").Value = -8333939
$ws.Range("N91").Value = -5365.4

# Row 97: Ore for Me
$ws.Range("H97").Value = 853.2414
$ws.Range("I97").Value = 749.7778
$ws.Range("J97").Value = 2250
$ws.Range("K97").Value = 749.7778
$ws.Range("L97").Value = 2250
$ws.Range("M97").Value = -253.7778
$ws.Range("N97").Value = -3242

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2265.1667
$ws.Range("I132").Value = 1718.8611
$ws.Range("K132").Value = 5156.5833
$ws.Range("M132").Value = -2626.5833

$ws = $wb.Worksheets.Item("BSM")
# Row 8: Stainless Chef
$ws.Range("H8").Value = 1842.1111
$ws.Range("I8").Value = 683.5
$ws.Range("J8").Value = 11111
$ws.Range("K8").Value = 683.5
$ws.Range("L8").Value = 11111
$ws.Range("M8").Value = -543.5
$ws.Range("N8").Value = -11391

# Row 10: Bring Me the Head Knife of Al'bedo Derssia
$ws.Range("H10").Value = 3679.6667
$ws.Range("I10").Value = 2993.3333
$ws.Range("K10").Value = 2993.3333
$ws.Range("M10").Value = -2853.3333

# Row 14: Farriers of Fortune
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# Row 15: Anutha Spatha
$ws.Range("H15").Value = 33000
$ws.Range("J15").Value = 33000
$ws.Range("L15").Value = 33000
$ws.Range("N15").Value = -33454

# Row 46: Spice Cadet
$ws.Range("H46").Value = 55555
$ws.Range("J46").Value = 55555
$ws.Range("L46").Value = 55555
$ws.Range("N46").Value = -56151

# Row 59: Pop That Top
$ws.Range("H59").Value = 118844
$ws.Range("J59").Value = 118844
$ws.Range("L59").Value = 118844
$ws.Range("N59").Value = -120538

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2554.1836
$ws.Range("I134").Value = 1450.2778
$ws.Range("J134").Value = 5611.154
$ws.Range("K134").Value = 4350.8334
$ws.Range("L134").Value = 16833.462
$ws.Range("M134").Value = -1815.8334
$ws.Range("N134").Value = -21903.462

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 35719816
$ws.Range("I31").Value = 2605.75
$ws.Range("J31").Value = 50006700
$ws.Range("K31").Value = 2605.75
$ws.Range("L31").Value = 50006700
$ws.Range("M31").Value = -2310.75
$ws.Range("N31").Value = -50007290

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 35719816
$ws.Range("I34").Value = 2605.75
$ws.Range("J34").Value = 50006700
$ws.Range("K34").Value = 2605.75
$ws.Range("L34").Value = 50006700
$ws.Range("M34").Value = -2403.75
$ws.Range("N34").Value = -50007104

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 2191.65
$ws.Range("I58").Value = 1804.22
$ws.Range("J58").Value = 4128.8
$ws.Range("K58").Value = 1804.22
$ws.Range("L58").Value = 4128.8
$ws.Range("M58").Value = -1601.22
$ws.Range("N58").Value = -4534.8

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 3656.6667
$ws.Range("I132").Value = 3508.7058
$ws.Range("J132").Value = 3908.2
$ws.Range("K132").Value = 10526.1174
$ws.Range("L132").Value = 11724.6
$ws.Range("M132").Value = -7996.117400000001
$ws.Range("N132").Value = -16784.6

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 8486.111000000001
$ws.Range("I134").Value = 21007.2
$ws.Range("J134").Value = 3670.3076
$ws.Range("K134").Value = 63021.60000000001
$ws.Range("L134").Value = 11010.9228
$ws.Range("M134").Value = -60486.60000000001
$ws.Range("N134").Value = -16080.9228

# Row 136: Turali Quality
$ws.Range("H136").Value = 2191.65
$ws.Range("I136").Value = 1804.22
$ws.Range("J136").Value = 4128.8
$ws.Range("K136").Value = 5412.66
$ws.Range("L136").Value = 12386.4
$ws.Range("M136").Value = -2862.66
$ws.Range("N136").Value = -17486.4

$ws = $wb.Worksheets.Item("CUL")
# Row 69: Loving That Muffin Top
$ws.Range("H69").Value = 1782
$ws.Range("I69").Value = 1033.3334
$ws.Range("J69").Value = 2905
$ws.Range("K69").Value = 3100.0002
$ws.Range("L69").Value = 8715
$ws.Range("M69").Value = -2289.0002
$ws.Range("N69").Value = -10337

# Row 72: Muffin of the Morn (L)
$ws.Range("H72").Value = 1782
$ws.Range("I72").Value = 1033.3334
$ws.Range("J72").Value = 2905
$ws.Range("K72").Value = 9300.000599999999
$ws.Range("L72").Value = 26145
$ws.Range("M72").Value = -5244.000599999999
$ws.Range("N72").Value = -34257

# Row 87: Soup That Eats Like a Knight
$ws.Range("H87").Value = 1000
$ws.Range("I87").Value = 1000
$ws.Range("K87").Value = 3000
$ws.Range("M87").Value = -1752

# Row 90: Like Ma Used to Make (L)
$ws.Range("H90").Value = 1000
$ws.Range("I90").Value = 1000
$ws.Range("K90").Value = 9000
$ws.Range("M90").Value = -2760

# Row 114: One Last Meal
$ws.Range("H114").Value = 3169.2222
$ws.Range("I114").Value = 100
$ws.Range("J114").Value = 3349.7646
$ws.Range("K114").Value = 300
$ws.Range("L114").Value = 10049.2938
$ws.Range("M114").Value = 2954
$ws.Range("N114").Value = -16557.2938

$ws = $wb.Worksheets.Item("GSM")
# Row 5: Hora at Me
$ws.Range("H5").Value = 10992.75
$ws.Range("J5").Value = 10992.75
$ws.Range("L5").Value = 10992.75
$ws.Range("N5").Value = -11216.75

# Row 46: Burning the Midnight Oil
$ws.Range("H46").Value = 32303.834
$ws.Range("J46").Value = 32303.834
$ws.Range("L46").Value = 32303.834
$ws.Range("N46").Value = -32615.834

# Row 52: It's My Business to Know Things
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# Row 103: Ring in the New
$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344

# Row 132: On Board for Lar
$ws.Range("H132").Value = 3137.4285
$ws.Range("I132").Value = 1568.4375
$ws.Range("J132").Value = 5229.4165
$ws.Range("K132").Value = 4705.3125
$ws.Range("L132").Value = 15688.2495
$ws.Range("M132").Value = -2175.3125
$ws.Range("N132").Value = -20748.2495

$ws = $wb.Worksheets.Item("LTW")
# Row 122: Hell on Leather
$ws.Range("H122").Value = 4283.448
$ws.Range("I122").Value = 2778.7778
$ws.Range("J122").Value = 6745.636
$ws.Range("K122").Value = 8336.3334
$ws.Range("L122").Value = 20236.908
$ws.Range("M122").Value = -5886.3334
$ws.Range("N122").Value = -25136.908

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 3310.3
$ws.Range("I136").Value = 1724.8096
$ws.Range("J136").Value = 7009.778
$ws.Range("K136").Value = 5174.4288
$ws.Range("L136").Value = 21029.334
$ws.Range("M136").Value = -2624.4288
$ws.Range("N136").Value = -26129.334

$ws = $wb.Worksheets.Item("WVR")
# Row 51: After the Smock-down
$ws.Range("H51").Value = 17341.334
$ws.Range("J51").Value = 18809.6
$ws.Range("L51").Value = 18809.6
$ws.Range("N51").Value = -19829.6

# Row 52: Party Animals
$ws.Range("H52").Value = 33347334
$ws.Range("I52").Value = 50001000
$ws.Range("J52").Value = 40000
$ws.Range("K52").Value = 50001000
$ws.Range("L52").Value = 40000
$ws.Range("M52").Value = -50000774
$ws.Range("N52").Value = -40452
